$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 8): A8 = 0.5 (time value), B8 = "c"
$ws.Range("A8").Value = 0.5
$ws.Range("B8").Value = "c"

# Expand the autofilter to cover the new row and add the new filter value
$ws.AutoFilterMode = $false
$ws.Range("A1:B8").AutoFilter(1, @("0.046", "0.500", "0.516"), 7) | Out-Null

# Keep the _FilterDatabase defined name in sync with the new autofilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Munka1!_FilterDatabase") {
        $n.RefersTo = "=Munka1!`$A`$1:`$B`$8"
    }
}

# Move the selection like in the edited file
$ws.Range("C7").Select() | Out-Null
